$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44826
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112036
$ws.Cells.Item($row, 7).Value = "Caigua"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 50
$ws.Cells.Item($row, 11).Value = 20000
$ws.Cells.Item($row, 12).Value = 20000
$ws.Cells.Item($row, 13).Value = 20000
$ws.Cells.Item($row, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 1333
$ws.Cells.Item($row, 17).Value = 15
$ws.Cells.Item($row, 18).Value = "Hortaliza"
